$d = $word.ActiveDocument

# The bullet reads "...directed fencing bouts at various fencing tournaments."
# Remove the redundant "fencing " right before "bouts" so it reads
# "...directed bouts at various fencing tournaments." Using a targeted
# Range().Delete() (instead of a Find/Replace spanning several runs) keeps
# the untouched surrounding runs - and their individual run formatting -
# intact, only removing the extra word and its leading space.

$anchor = $d.Content
$found = $anchor.Find.Execute("directed fencing bouts")

if ($found) {
    $wordStart = $anchor.Start + "directed".Length
    $wordEnd = $wordStart + " fencing".Length
    $d.Range($wordStart, $wordEnd).Delete()
}
